$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Typ" column
$ws.Range("B1").Value = "Typ"

# Existing MCD addresses (rows 2-9) get a Type of "MCD"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = "MCD"
}

# Header for new "Priorytety" column
$ws.Range("C1").Value = "Priorytety"

# Priority 1 for the existing MCD addresses
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}

# New KFC addresses (rows 10-12)
$ws.Range("A10").Value = "al. Wincentego Witosa 32, 20-315 Lublin"
$ws.Range("A11").Value = "al. Tysiąclecia 12, 20-121 Lublin"
$ws.Range("A12").Value = "Ireny Sendlerowej 1, 20-817 Lublin"

for ($r = 10; $r -le 12; $r++) {
    $ws.Cells.Item($r, 2).Value = "KFC"
    $ws.Cells.Item($r, 3).Value = 2
}

# Update the active selection to match the target state
$ws.Range("E10").Select()
